$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cells (e.g. H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the new I and J columns of data for rows 2-13.
$ws.Cells.Item(2, 9).Value = 7
$ws.Cells.Item(2, 10).Value = 9

$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 6

$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 5

$ws.Cells.Item(5, 9).Value = 1
$ws.Cells.Item(5, 10).Value = 2

$ws.Cells.Item(6, 9).Value = 6
$ws.Cells.Item(6, 10).Value = 7

$ws.Cells.Item(7, 9).Value = 6
$ws.Cells.Item(7, 10).Value = 7

$ws.Cells.Item(8, 9).Value = 6
$ws.Cells.Item(8, 10).Value = 7

$ws.Cells.Item(9, 9).Value = 9
$ws.Cells.Item(9, 10).Value = 9

$ws.Cells.Item(10, 9).Value = 6
$ws.Cells.Item(10, 10).Value = 6

$ws.Cells.Item(11, 9).Value = 6
$ws.Cells.Item(11, 10).Value = 7

$ws.Cells.Item(12, 9).Value = 6
$ws.Cells.Item(12, 10).Value = 7

$ws.Cells.Item(13, 9).Value = 1
$ws.Cells.Item(13, 10).Value = 2
